$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.349.73"
$ws.Range("E2").Value = "  -1.85%  "
$ws.Range("D3").Value = "1.656.15"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.29"
$ws.Range("E5").Value = "  -0.69%  "
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.67"
$ws.Range("E8").Value = "  +0.45%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("E10").Value = "  -1.12%  "
$ws.Range("E11").Value = "  -0.39%  "
$ws.Range("D12").Value = "1.890.04"
$ws.Range("E12").Value = "  -0.78%  "
$ws.Range("D13").Value = "1.655.60"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.572"
$ws.Range("E14").Value = "  +3.69%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.08"
$ws.Range("E15").Value = "  -1.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.81"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("D17").Value = "27.350.31"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "232.00"
$ws.Range("E18").Value = "  -7.79%  "
$ws.Range("E19").Value = "  -0.86%  "
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.38"
$ws.Range("E22").Value = "  -2.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.37"
$ws.Range("E23").Value = "  +0.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "15.90"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("E28").Value = "  -0.18%  "
$ws.Range("E29").Value = "  -0.43%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  -0.84%  "
$ws.Range("E31").Value = "  -3.88%  "
$ws.Range("E32").Value = "  -1.33%  "
$ws.Range("D33").Value = "1.448.36"
$ws.Range("E33").Value = "  +1.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.13"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("E36").Value = "  -0.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.908"
$ws.Range("E37").Value = "  -2.31%  "
$ws.Range("E38").Value = "  -2.07%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.53"
$ws.Range("E42").Value = "  +2.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.15"
$ws.Range("E43").Value = "  -6.66%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("D45").Value = "1.797.88"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.70"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.13"
$ws.Range("D49").Value = "0.0₆0106"
$ws.Range("E49").Value = "  -4.85%  "
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.72"
$ws.Range("E51").Value = "  -1.67%  "
